$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.791.96"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "1.625.49"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5106"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2559"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06318"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.37"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07770"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.221"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "1.625.92"
$ws.Range("E13").Value = "  -0.84%  "
$ws.Range("D14").Value = "1.846.37"
$ws.Range("E14").Value = "  -0.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5526"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.08%  "
$ws.Range("E16").Value = "  -0.99%  "
$ws.Range("D17").Value = "0.0₅7499"
$ws.Range("E17").Value = "  -2.48%  "
$ws.Range("D18").Value = "25.790.82"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.405"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.753"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.993"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.869"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1241"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.83%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.699"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.71%  "
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04851"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.234"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.157"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.534"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.364"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.8912"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.536"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5497"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.91%  "
$ws.Range("D39").Value = "1.114.06"
$ws.Range("E39").Value = "  -2.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01546"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.19%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.509"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7956"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "96.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.16%  "
$ws.Range("D45").Value = "1.770.37"
$ws.Range("E45").Value = "  -0.34%  "
$ws.Range("D46").Value = "0.0₈116"
$ws.Range("E46").Value = "  -6.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4424"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9971"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.48"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05123"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.523"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.70%  "
